$p = $ppt.ActivePresentation

# --- Slide 2 ("Sections to Cover") ---
# "Directives" bullet -> "Modules" text + hyperlink to the Modules slide (slide 5)
# "Controllers" bullet -> keep text, add hyperlink to the Controllers slide (slide 6)
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$paraDirectives = $tr.Paragraphs(2)
$charsDirectives = $paraDirectives.Characters(1, $paraDirectives.Length)
$charsDirectives.Text = "Modules"
$actionDirectives = $charsDirectives.ActionSettings(1)
$actionDirectives.Action = 7
$hlDirectives = $actionDirectives.Hyperlink
$hlDirectives.Address = ""
$hlDirectives.SubAddress = "5,5,Modules"

$paraControllers = $tr.Paragraphs(3)
$charsControllers = $paraControllers.Characters(1, $paraControllers.Length)
$actionControllers = $charsControllers.ActionSettings(1)
$actionControllers.Action = 7
$hlControllers = $actionControllers.Hyperlink
$hlControllers.Address = ""
$hlControllers.SubAddress = "6,6,Controllers"

# --- Slide 5 ("Introduction - Modules") ---
# Title shape: drop the leading "Introduction" + line break + "  - " prefix,
# leaving just the "Modules" run (sz=2800) in place.
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Item(1)
$titleRange = $title5.TextFrame.TextRange
$prefixLen = $titleRange.Length - 7
$prefix = $titleRange.Characters(1, $prefixLen)
$prefix.Text = ""
